# Automatische test-sync: 2025-06-26 19:07:50
# Append the new "Van Ommeren" e-mail as row 7 on the Logs sheet, extend the
# conditional-formatting ranges to cover it, and refresh the Dashboard
# category-summary table (Bestelling/Levering now sorts to 2nd place with
# count 2; Retour/Terugbetaling and Productinformatie shift down one row).

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# ---------------------------------------------------------------------
# 1. New row 7 on the Logs sheet
# ---------------------------------------------------------------------
$logs.Range("A7").Value = "Kan je voor mij 2500 M8 slotbouten bestellen bij Van Ommeren?"
$logs.Range("B7").Value = "MailMind Test <mailmind.test@zohomail.eu>"

$logs.Range("C7").Value = @"
Hoi, 
Kun jij 2500 stuks M8 slotbouten bestellen bij Van Ommeren voor levering op donderdag? 
Laat me weten als dit gelukt is.
Groeten, 
Jeroen
Sent using {0}
"@

$logs.Range("D7").Value = "Bestelling / Levering"

$logs.Range("E7").Value = @"
Beste Jeroen,
Bedankt voor je bericht. Ik ben een geautomatiseerde assistent en kan geen bestellingen plaatsen. Je kunt rechtstreeks contact opnemen met Van Ommeren om je bestelling van 2500 stuks M8 slotbouten te plaatsen voor levering op donderdag.
Voor verdere assistentie kun je mij altijd contacteren.
Met vriendelijke groet,
[Bedrijfsnaam] E-mailassistent
"@

$logs.Range("F7").Value = "2025-06-26 19:07:34"
$logs.Range("G7").Value = "Ja"
$logs.Range("H7").Value = "Nee"

# ---------------------------------------------------------------------
# 2. Extend the three conditional-formatting blocks from row 6 to row 7
# ---------------------------------------------------------------------
$logs.Range("D2:D6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D7"))
$logs.Range("G2:G6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G7"))
$logs.Range("H2:H6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H7"))

# ---------------------------------------------------------------------
# 3. Refresh the Dashboard category-summary table
#    (row 2 "Offerte / Prijsaanvraag" / 2 is unchanged)
# ---------------------------------------------------------------------
$dashboard.Range("A3").Value = "Bestelling / Levering"
$dashboard.Range("B3").Value = 2
$dashboard.Range("A4").Value = "Retour / Terugbetaling"
$dashboard.Range("B4").Value = 1
$dashboard.Range("A5").Value = "Productinformatie"
$dashboard.Range("B5").Value = 1
